# p2p_parser: show bank transfers in the final result in a single column
#
# The three result sheets each had separate "Einzahlungen" (deposits) and
# "Auszahlungen" (withdrawals) columns. They are merged into a single
# "Ein-/Auszahlungen" column: rename the "Einzahlungen" column header and
# delete the "Auszahlungen" column that immediately follows it.

$wb = $excel.ActiveWorkbook

# --- Monatsergebnisse: rename + delete its "Auszahlungen" column first,
# and leave its selection on the cell the sheet ends up showing.
$ws2 = $wb.Worksheets.Item("Monatsergebnisse")
$ws2.Range("F1").Value = "Ein-/Auszahlungen"
$ws2.Range("G1").EntireColumn.Delete()
$ws2.Range("G1").Select()

# --- Gesamtergebnis: same merge, one column earlier since it has no
# Datum/Monat column.
$ws3 = $wb.Worksheets.Item("Gesamtergebnis")
$ws3.Range("E1").Value = "Ein-/Auszahlungen"
$ws3.Range("F1").EntireColumn.Delete()
$ws3.Range("F1").Select()

# --- Tagesergebnisse: same merge, then make it the active sheet/tab.
$ws1 = $wb.Worksheets.Item("Tagesergebnisse")
$ws1.Range("F1").Value = "Ein-/Auszahlungen"
$ws1.Range("G1").EntireColumn.Delete()
$ws1.Activate()
$ws1.Range("F7").Select()
